# Apply the data correction: swap the values in columns C and D for the
# affected rows (the original upload had the "left" / "right" choice
# counts transposed for these observations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5,7,8,9,10,11,12,15,16,17,18,19,20,21,22,24,29,32,33,40,44,45,47,48,49,51,52,53,54,55,56,57,60,61,63,64,66,67,68,70,71,72,73,74,75,76)

foreach ($r in $rows) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}

# Update the view state to reflect the saved workbook (zoom level and
# selected cell on Sheet1).
$ws.Range("G9").Select()
$excel.ActiveWindow.Zoom = 62
